$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top (using original row numbers) so earlier
# deletions don't shift the row numbers of rows we still need to delete.

# Delete ISABEL row (005624730, 5970.68) - originally row 15
$ws.Range("A15").EntireRow.Delete() | Out-Null

# Delete MARIA row (001090818, 19646.89) - originally row 7
$ws.Range("A7").EntireRow.Delete() | Out-Null

# Delete RAPHAELA row (005366255, 30001.47) - originally row 5
$ws.Range("A5").EntireRow.Delete() | Out-Null

# Delete CINTIA row (005949170, 166000) - originally row 4
$ws.Range("A4").EntireRow.Delete() | Out-Null

# Delete CARLOS row (005696533, 202032.82) - originally row 3
$ws.Range("A3").EntireRow.Delete() | Out-Null

# Now LEDA (002636063) has shifted from original row 6 to row 3.
# Update her balance from 22556.07 to 33000.
$ws.Cells.Item(3, 3).Value = 33000

# DOUGLAS (originally row 21, now row 16) account number changes
# from 004384167 to 005608744. Force the value to be stored as text so
# the leading zeros are preserved, then clear the number format so the
# cell keeps the default (unstyled) look it had before.
$cell = $ws.Cells.Item(16, 1)
$cell.NumberFormat = "@"
$cell.Value = "005608744"
$cell.ClearFormats()
